$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their Text format so Excel does not
# re-interpret the numeric-looking / percent-looking strings as
# numbers when the new values are written.
$targets = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5",
    "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9",
    "D10", "E10", "D11", "E11", "E12", "D13", "E13", "D14",
    "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18",
    "E18", "D19", "E19", "D20", "E20", "E21", "E22", "D23",
    "E23", "E24", "D25", "E25", "D26", "E26", "E27", "E28",
    "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44",
    "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48",
    "D49", "E49", "D50", "E50"
)
foreach ($addr in $targets) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed price / 1h-volume values scraped on
# Tue Jan  3 08:32:04 UTC 2023.
$ws.Range("D2").Value = "246.05"
$ws.Range("E2").Value = "-0.14%"
$ws.Range("D3").Value = "30.01"
$ws.Range("E3").Value = "0.31%"
$ws.Range("D4").Value = "5.154"
$ws.Range("E4").Value = "-0.19%"
$ws.Range("D5").Value = "0.05778"
$ws.Range("E5").Value = "0.86%"
$ws.Range("D6").Value = "6.653"
$ws.Range("E6").Value = "1.06%"
$ws.Range("D7").Value = "3.233"
$ws.Range("E7").Value = "6.71%"
$ws.Range("D8").Value = "0.8494"
$ws.Range("E8").Value = "-0.82%"
$ws.Range("D9").Value = "0.8556"
$ws.Range("E9").Value = "-2.02%"
$ws.Range("D10").Value = "0.1381"
$ws.Range("E10").Value = "1.32%"
$ws.Range("D11").Value = "0.07096"
$ws.Range("E11").Value = "1.53%"
$ws.Range("E12").Value = "11.94%"
$ws.Range("D13").Value = "0.09378"
$ws.Range("E13").Value = "0.02%"
$ws.Range("D14").Value = "0.001525"
$ws.Range("E14").Value = "0.64%"
$ws.Range("D15").Value = "0.0005982"
$ws.Range("E15").Value = "-0.58%"
$ws.Range("D16").Value = "0.006077"
$ws.Range("E16").Value = "0.57%"
$ws.Range("D17").Value = "3.508"
$ws.Range("E17").Value = "-0.06%"
$ws.Range("D18").Value = "2.222"
$ws.Range("E18").Value = "2.19%"
$ws.Range("D19").Value = "0.3160"
$ws.Range("E19").Value = "0.52%"
$ws.Range("D20").Value = "0.03374"
$ws.Range("E20").Value = "1.75%"
$ws.Range("E21").Value = "-0.61%"
$ws.Range("E22").Value = "-3.45%"
$ws.Range("D23").Value = "0.04121"
$ws.Range("E23").Value = "-0.90%"
$ws.Range("E24").Value = "0.30%"
$ws.Range("D25").Value = "0.001226"
$ws.Range("E25").Value = "1.27%"
$ws.Range("D26").Value = "0.004142"
$ws.Range("E26").Value = "-7.92%"
$ws.Range("E27").Value = "1.88%"
$ws.Range("E28").Value = "5.26%"
$ws.Range("D40").Value = "0.03754"
$ws.Range("E40").Value = "-0.87%"
$ws.Range("D41").Value = "0.1069"
$ws.Range("E41").Value = "0.15%"
$ws.Range("E42").Value = "-4.46%"
$ws.Range("D43").Value = "0.003526"
$ws.Range("E43").Value = "-38.89%"
$ws.Range("D44").Value = "0.008836"
$ws.Range("E44").Value = "-11.74%"
$ws.Range("D45").Value = "0.00005430"
$ws.Range("E45").Value = "6.57%"
$ws.Range("E46").Value = "0.26%"
$ws.Range("D47").Value = "0.07102"
$ws.Range("E47").Value = "-11.02%"
$ws.Range("D48").Value = "0.002190"
$ws.Range("E48").Value = "-19.53%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.26%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.26%"
